$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits after the
#    "...' family" run (end of the "fonts" list item).
# ---------------------------------------------------------------------
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # no pre-existing bookmark - nothing to remove
}

# ---------------------------------------------------------------------
# 2) Split the run "Introducing let & const variables" into two runs:
#      "Introducing let & const "  +  " for variable declaration"
#    (same run formatting / Times New Roman 14pt for both).
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Introducing let & const variables", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Shrink the found range down to just the first half of the text ...
$r.Text = "Introducing let & const "
$r.Collapse(0)

# ... then append the second half plus a throwaway marker run. Keeping a
# trailing marker run (instead of ending the paragraph right on the new
# text) avoids a boundary quirk where a bookmark collapsed exactly at the
# end of a paragraph's last run gets mis-anchored back at the paragraph's
# start; the marker is removed again below once the bookmark is safely
# placed.
$r.InsertAfter(" for variable declarationGOBACKMARKER")
$r.Font.Name = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

# ---------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark right after the new text (i.e. at
#    the end of this paragraph, matching where Word leaves it following
#    the most recent edit).
# ---------------------------------------------------------------------
$markerRange = $d.Content
$markerRange.Find.Execute("GOBACKMARKER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markerRange)

# Remove the throwaway marker text now that the bookmark has a stable
# anchor.
$cleanupRange = $d.Content
$cleanupRange.Find.Execute("GOBACKMARKER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cleanupRange.Text = ""
